# "lookup broker SSI built" - A5 now holds the numeric trade/ticket id (was
# text "19437-A"), and the settle/as-of dates in rows 6-7 (previously plain
# text like "15/10/20"/"21/10/20") become real date values so a broker SSI
# lookup can use them like the already-numeric dates in rows 4-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A5: text "19437-A" -> number 40019
$ws.Range("A5").Value = 40019

# Reuse the existing date number-format/style already applied to H4:I4 (the
# "As of Dt"/"Stl Date" columns in the rows above) so H6:I7 end up sharing
# the very same style instead of minting a new one.
$ws.Range("H4:I4").Copy()
$ws.Range("H6:I6").PasteSpecial(-4122)
$ws.Range("H7:I7").PasteSpecial(-4122)

# Write the real date values (as date serials) into H6/I6/H7/I7.
$ws.Range("H6").Value = 44119
$ws.Range("I6").Value = 44125
$ws.Range("H7").Value = 44119
$ws.Range("I7").Value = 44125

$ws.Range("J11").Select()
